# Update the hourly regression-output table (cap_gen_year16final) with
# refreshed coefficients / stats from the latest model run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 0.1053331483093296

$ws.Cells.Item(3, 2).Value = 0.1136322875517179
$ws.Cells.Item(3, 8).Value = 0.2189654358610475

$ws.Cells.Item(4, 2).Value = 0.1024582439816771
$ws.Cells.Item(4, 8).Value = 0.2077913922910067

$ws.Cells.Item(5, 2).Value = 0.04180458250060327
$ws.Cells.Item(5, 8).Value = 0.1471377308099329

$ws.Cells.Item(6, 2).Value = 0.02740139629355406
$ws.Cells.Item(6, 3).Value = 0.002515299710752093
$ws.Cells.Item(6, 4).Value = 3.993518845385184
$ws.Cells.Item(6, 5).Value = 0.01811784907453557
$ws.Cells.Item(6, 6).Value = 0.02246504313777602
$ws.Cells.Item(6, 7).Value = 0.03233774944933174
$ws.Cells.Item(6, 8).Value = 0.1327345446028836

$ws.Cells.Item(7, 2).Value = 0.01935738793585521
$ws.Cells.Item(7, 3).Value = 0.002275165851404394
$ws.Cells.Item(7, 4).Value = 2.363449289881279
$ws.Cells.Item(7, 5).Value = 0.04304667101166117
$ws.Cells.Item(7, 6).Value = 0.01489348645498139
$ws.Cells.Item(7, 7).Value = 0.02382128941672922
$ws.Cells.Item(7, 8).Value = 0.1246905362451848

$ws.Cells.Item(8, 2).Value = 0.01859274452294621
$ws.Cells.Item(8, 3).Value = 0.002034456734881975
$ws.Cells.Item(8, 4).Value = 2.305548800137076
$ws.Cells.Item(8, 5).Value = 0.0172915928024556
$ws.Cells.Item(8, 6).Value = 0.01459999903187738
$ws.Cells.Item(8, 7).Value = 0.02258549001401522
$ws.Cells.Item(8, 8).Value = 0.1239258928322758

$ws.Cells.Item(9, 2).Value = 0.01341835699473661
$ws.Cells.Item(9, 3).Value = 0.001917983871129994
$ws.Cells.Item(9, 4).Value = 1.958452738426621
$ws.Cells.Item(9, 5).Value = 0.00580895883881699
$ws.Cells.Item(9, 6).Value = 0.009653782993790718
$ws.Cells.Item(9, 7).Value = 0.0171829309956823
$ws.Cells.Item(9, 8).Value = 0.1187515053040662

$ws.Cells.Item(10, 2).Value = 0.01500699362001514
$ws.Cells.Item(10, 3).Value = 0.002552263456055266
$ws.Cells.Item(10, 4).Value = 2.12212481948439
$ws.Cells.Item(10, 5).Value = 0.00452888606977613
$ws.Cells.Item(10, 6).Value = 0.009996370829102581
$ws.Cells.Item(10, 7).Value = 0.02001761641092765
$ws.Cells.Item(10, 8).Value = 0.1203401419293447

$ws.Cells.Item(11, 2).Value = 0.02943687797160316
$ws.Cells.Item(11, 8).Value = 0.1347700262809327

$ws.Cells.Item(12, 2).Value = 0.0540681740750615
$ws.Cells.Item(12, 8).Value = 0.1594013223843911

$ws.Cells.Item(13, 2).Value = 0.06931410665127544
$ws.Cells.Item(13, 8).Value = 0.174647254960605

$ws.Cells.Item(14, 2).Value = 0.07694405753430401
$ws.Cells.Item(14, 8).Value = 0.1822772058436336

$ws.Cells.Item(15, 2).Value = 0.08434279519203611
$ws.Cells.Item(15, 8).Value = 0.1896759435013657

$ws.Cells.Item(16, 2).Value = 0.08793985650266503
$ws.Cells.Item(16, 8).Value = 0.1932730048119946

$ws.Cells.Item(17, 2).Value = 0.09051673174416253
$ws.Cells.Item(17, 8).Value = 0.1958498800534921

$ws.Cells.Item(18, 2).Value = -0.1053331483093296

$ws.Cells.Item(19, 2).Value = 0.09212893575831305
$ws.Cells.Item(19, 8).Value = 0.1974620840676426

$ws.Cells.Item(20, 2).Value = 0.09550757932538084
$ws.Cells.Item(20, 8).Value = 0.2008407276347104

$ws.Cells.Item(21, 2).Value = 0.1001145206728902
$ws.Cells.Item(21, 8).Value = 0.2054476689822198

$ws.Cells.Item(22, 2).Value = 0.1041123371286635
$ws.Cells.Item(22, 3).Value = 0.006734902740321215
$ws.Cells.Item(22, 4).Value = 235140220804.5638
$ws.Cells.Item(22, 5).Value = 0.0321794941641154
$ws.Cells.Item(22, 6).Value = 0.09088387040358084
$ws.Cells.Item(22, 7).Value = 0.1173408038537467
$ws.Cells.Item(22, 8).Value = 0.2094454854379931

$ws.Cells.Item(23, 2).Value = 0.1084776228961321
$ws.Cells.Item(23, 3).Value = 0.006985265189092959
$ws.Cells.Item(23, 4).Value = 359274961006.428
$ws.Cells.Item(23, 5).Value = 0.03653729407134387
$ws.Cells.Item(23, 6).Value = 0.09475932487961147
$ws.Cells.Item(23, 7).Value = 0.1221959209126524
$ws.Cells.Item(23, 8).Value = 0.2138107712054617

$ws.Cells.Item(24, 2).Value = 0.1114466841574084
$ws.Cells.Item(24, 8).Value = 0.2167798324667379

$ws.Cells.Item(25, 2).Value = 0.1130090738038105
$ws.Cells.Item(25, 3).Value = 0.007274955822356711
$ws.Cells.Item(25, 4).Value = 26.89075381159514
$ws.Cells.Item(25, 5).Value = 0.04669545581823498
$ws.Cells.Item(25, 6).Value = 0.09870644857410409
$ws.Cells.Item(25, 7).Value = 0.1273116990335168
$ws.Cells.Item(25, 8).Value = 0.21834222211314

$ws.Cells.Item(26, 2).Value = 0.1158079375841897
$ws.Cells.Item(26, 8).Value = 0.2211410858935193

$ws.Cells.Item(27, 2).Value = 0.1209846171645378
$ws.Cells.Item(27, 3).Value = 0.006979416414990986
$ws.Cells.Item(27, 4).Value = 26.73276823934978
$ws.Cells.Item(27, 5).Value = 0.05040596834229742
$ws.Cells.Item(27, 6).Value = 0.1072698523614101
$ws.Cells.Item(27, 7).Value = 0.1346993819676664
$ws.Cells.Item(27, 8).Value = 0.2263177654738673

$ws.Cells.Item(28, 2).Value = 0.1224983424820655
$ws.Cells.Item(28, 3).Value = 0.007210199647030541
$ws.Cells.Item(28, 4).Value = 25.46995059647357
$ws.Cells.Item(28, 5).Value = 0.08297233238895797
$ws.Cells.Item(28, 6).Value = 0.108347778023142
$ws.Cells.Item(28, 7).Value = 0.1366489069409896
$ws.Cells.Item(28, 8).Value = 0.227831490791395

$ws.Cells.Item(29, 2).Value = 0.01972539771923136
$ws.Cells.Item(29, 3).Value = 0.001774406468956909
$ws.Cells.Item(29, 4).Value = 2.822369345998223
$ws.Cells.Item(29, 5).Value = 0.06167796232882037
$ws.Cells.Item(29, 6).Value = 0.0162356048534354
$ws.Cells.Item(29, 7).Value = 0.02321519058502691
$ws.Cells.Item(29, 8).Value = 0.125058546028561
